$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 5801.75
$ws.Range("I11").Value = 5801.75
$ws.Range("K11").Value = 5801.75
$ws.Range("M11").Value = -5661.75
$ws.Range("H18").Value = 1798.5
$ws.Range("I18").Value = 1798.5
$ws.Range("K18").Value = 1798.5
$ws.Range("M18").Value = -1514.5
$ws.Range("H40").Value = 5700.2856
$ws.Range("J40").Value = 5817
$ws.Range("L40").Value = 5817
$ws.Range("N40").Value = -6167
$ws.Range("H68").Value = 40247.5
$ws.Range("J68").Value = 40247.5
$ws.Range("L68").Value = 40247.5
$ws.Range("N68").Value = -41745.5
$ws.Range("H69").Value = 15139.5
$ws.Range("I69").Value = 7336
$ws.Range("K69").Value = 22008
$ws.Range("M69").Value = -21134
$ws.Range("H71").Value = 40247.5
$ws.Range("J71").Value = 40247.5
$ws.Range("L71").Value = 120742.5
$ws.Range("N71").Value = -128230.5
$ws.Range("H72").Value = 15139.5
$ws.Range("I72").Value = 7336
$ws.Range("K72").Value = 66024
$ws.Range("M72").Value = -61656
$ws.Range("H80").Value = 2549.7144
$ws.Range("J80").Value = 3299.5
$ws.Range("L80").Value = 9898.5
$ws.Range("N80").Value = -11894.5
$ws.Range("H83").Value = 2549.7144
$ws.Range("J83").Value = 3299.5
$ws.Range("L83").Value = 29695.5
$ws.Range("N83").Value = -39679.5
$ws.Range("H92").Value = 5982.5835
$ws.Range("I92").Value = 6864.3
$ws.Range("J92").Value = 1574
$ws.Range("K92").Value = 6864.3
$ws.Range("L92").Value = 1574
$ws.Range("M92").Value = -5616.3
$ws.Range("N92").Value = -4070
$ws.Range("H113").Value = 14030
$ws.Range("I113").Value = 17322.5
$ws.Range("K113").Value = 17322.5
$ws.Range("M113").Value = -14068.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8136.875
$ws.Range("J2").Value = 8985
$ws.Range("L2").Value = 8985
$ws.Range("N2").Value = -9211
$ws.Range("H5").Value = 2120089
$ws.Range("I5").Value = 2609294.2
$ws.Range("K5").Value = 2609294.2
$ws.Range("M5").Value = -2609182.2
$ws.Range("H102").Value = 22016.584
$ws.Range("I102").Value = 2366.6667
$ws.Range("K102").Value = 2366.6667
$ws.Range("M102").Value = -744.6667000000002
$ws.Range("H110").Value = 4546442.5
$ws.Range("I110").Value = 4785676
$ws.Range("K110").Value = 4785676
$ws.Range("M110").Value = -4783631
$ws.Range("H116").Value = 8136.875
$ws.Range("J116").Value = 8985
$ws.Range("L116").Value = 8985
$ws.Range("N116").Value = -13573

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8136.875
$ws.Range("J3").Value = 8985
$ws.Range("L3").Value = 8985
$ws.Range("N3").Value = -9213
$ws.Range("H4").Value = 2120089
$ws.Range("I4").Value = 2609294.2
$ws.Range("K4").Value = 2609294.2
$ws.Range("M4").Value = -2609179.2
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H99").Value = 7474.8335
$ws.Range("I99").Value = 1370
$ws.Range("K99").Value = 1370
$ws.Range("M99").Value = 128

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 730.4706
$ws.Range("I7").Value = 95.833336
$ws.Range("K7").Value = 95.833336
$ws.Range("M7").Value = 17.166664
$ws.Range("H16").Value = 10998.8
$ws.Range("I16").Value = 1248.75
$ws.Range("K16").Value = 1248.75
$ws.Range("M16").Value = -961.75
$ws.Range("H86").Value = 18150.75
$ws.Range("I86").Value = 21689.334
$ws.Range("J86").Value = 7535
$ws.Range("K86").Value = 21689.334
$ws.Range("L86").Value = 7535
$ws.Range("M86").Value = -20566.334
$ws.Range("N86").Value = -9781
$ws.Range("H89").Value = 18150.75
$ws.Range("I89").Value = 21689.334
$ws.Range("J89").Value = 7535
$ws.Range("K89").Value = 108446.67
$ws.Range("L89").Value = 37675
$ws.Range("M89").Value = -102830.67
$ws.Range("N89").Value = -48907
$ws.Range("H93").Value = 17656.625
$ws.Range("I93").Value = 17393.285
$ws.Range("J93").Value = 19500
$ws.Range("K93").Value = 17393.285
$ws.Range("L93").Value = 19500
$ws.Range("M93").Value = -15521.285
$ws.Range("N93").Value = -23244
$ws.Range("H99").Value = 3032954.5
$ws.Range("I99").Value = 3959406
$ws.Range("J99").Value = 2511825.8
$ws.Range("K99").Value = 3959406
$ws.Range("L99").Value = 2511825.8
$ws.Range("M99").Value = -3957908
$ws.Range("N99").Value = -2514821.8
$ws.Range("H107").Value = 2695.516
$ws.Range("I107").Value = 1250.5264
$ws.Range("J107").Value = 4983.4165
$ws.Range("K107").Value = 1250.5264
$ws.Range("L107").Value = 4983.4165
$ws.Range("M107").Value = 669.4736
$ws.Range("N107").Value = -8823.416499999999
$ws.Range("H113").Value = 10998.8
$ws.Range("I113").Value = 1248.75
$ws.Range("K113").Value = 1248.75
$ws.Range("M113").Value = 921.25
$ws.Range("H126").Value = 3032954.5
$ws.Range("I126").Value = 3959406
$ws.Range("J126").Value = 2511825.8
$ws.Range("K126").Value = 11878218
$ws.Range("L126").Value = 7535477.399999999
$ws.Range("M126").Value = -11875748
$ws.Range("N126").Value = -7540417.399999999
$ws.Range("H134").Value = 35721900
$ws.Range("J134").Value = 76937220
$ws.Range("L134").Value = 230811660
$ws.Range("N134").Value = -230816730

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 484.85
$ws.Range("I40").Value = 543
$ws.Range("J40").Value = 252.25
$ws.Range("K40").Value = 2172
$ws.Range("L40").Value = 1009
$ws.Range("M40").Value = -2103
$ws.Range("N40").Value = -1147
$ws.Range("H136").Value = 2468
$ws.Range("I136").Value = 1727
$ws.Range("K136").Value = 5181
$ws.Range("M136").Value = -81
$ws.Range("H137").Value = 1735.25
$ws.Range("J137").Value = 1822
$ws.Range("L137").Value = 5466
$ws.Range("N137").Value = -15666

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 29482.285
$ws.Range("I80").Value = 24194.777
$ws.Range("K80").Value = 24194.777
$ws.Range("M80").Value = -23196.777
$ws.Range("H83").Value = 29482.285
$ws.Range("I83").Value = 24194.777
$ws.Range("K83").Value = 120973.885
$ws.Range("M83").Value = -115981.885
$ws.Range("H97").Value = 6120.7
$ws.Range("I97").Value = 1622.5
$ws.Range("J97").Value = 16616.5
$ws.Range("K97").Value = 1622.5
$ws.Range("L97").Value = 16616.5
$ws.Range("M97").Value = -1126.5
$ws.Range("N97").Value = -17608.5
$ws.Range("H113").Value = 6227
$ws.Range("I113").Value = 6454
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 6454
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -4284
$ws.Range("N113").Value = -10340

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2107170.5
$ws.Range("I40").Value = 2811.9285
$ws.Range("J40").Value = 4211529
$ws.Range("K40").Value = 2811.9285
$ws.Range("L40").Value = 4211529
$ws.Range("M40").Value = -2675.9285
$ws.Range("N40").Value = -4211801
$ws.Range("H58").Value = 12689.667
$ws.Range("I58").Value = 9034.5
$ws.Range("K58").Value = 9034.5
$ws.Range("M58").Value = -8774.5
$ws.Range("H68").Value = 5487.1577
$ws.Range("I68").Value = 3215.2856
$ws.Range("J68").Value = 6812.4165
$ws.Range("K68").Value = 3215.2856
$ws.Range("L68").Value = 6812.4165
$ws.Range("M68").Value = -2466.2856
$ws.Range("N68").Value = -8310.416499999999
$ws.Range("H71").Value = 5487.1577
$ws.Range("I71").Value = 3215.2856
$ws.Range("J71").Value = 6812.4165
$ws.Range("K71").Value = 16076.428
$ws.Range("L71").Value = 34062.0825
$ws.Range("M71").Value = -12332.428
$ws.Range("N71").Value = -41550.0825
$ws.Range("H82").Value = 7060.353
$ws.Range("I82").Value = 4425
$ws.Range("J82").Value = 15625.25
$ws.Range("K82").Value = 4425
$ws.Range("L82").Value = 15625.25
$ws.Range("M82").Value = -4064
$ws.Range("N82").Value = -16347.25
$ws.Range("H85").Value = 7060.353
$ws.Range("I85").Value = 4425
$ws.Range("J85").Value = 15625.25
$ws.Range("K85").Value = 4425
$ws.Range("L85").Value = 15625.25
$ws.Range("M85").Value = -3177
$ws.Range("N85").Value = -18121.25
$ws.Range("H122").Value = 27992516
$ws.Range("I122").Value = 43398910
$ws.Range("J122").Value = 3342285
$ws.Range("K122").Value = 130196730
$ws.Range("L122").Value = 10026855
$ws.Range("M122").Value = -130194280
$ws.Range("N122").Value = -10031755
$ws.Range("H132").Value = 1754063.1
$ws.Range("I132").Value = 7557.364
$ws.Range("J132").Value = 3355026.8
$ws.Range("K132").Value = 22672.092
$ws.Range("L132").Value = 10065080.4
$ws.Range("M132").Value = -20142.092
$ws.Range("N132").Value = -10070140.4
$ws.Range("H136").Value = 12848.69
$ws.Range("I136").Value = 13151
$ws.Range("K136").Value = 39453
$ws.Range("M136").Value = -36903

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 456922.97
$ws.Range("I122").Value = 601972
$ws.Range("K122").Value = 1805916
$ws.Range("M122").Value = -1803466
$ws.Range("H126").Value = 6682693
$ws.Range("I126").Value = 7722.222
$ws.Range("K126").Value = 23166.666
$ws.Range("M126").Value = -20696.666
$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920
